$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").NumberFormat = "@"
$ws.Range("A20").Value = "05/06/2025"
$ws.Range("B20").Value = 550.3459999999977
$ws.Range("C20").Value = 0.09085193678158869
$ws.Range("D20").Value = 50
